$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# A2: drop the explicit cell style (s="1" -> default/no style)
$ws.Range("A2").Style = "Normal"

# B2: FARM_SITE_ID -> a UNIQUE_FLOCK_ID-style value that contains zero-width
# spaces interleaved with the digits. Built with -join so the interpreter
# keeps it as a literal string instead of mis-parsing it as a number.
$zwsp = [char]0x200B
$uniqueFlockId = ($zwsp, "1", $zwsp, "0", $zwsp, "0", $zwsp, "1", $zwsp, "1", $zwsp, "9", $zwsp, "0", $zwsp) -join ""
$ws.Range("B2").Value = $uniqueFlockId

# L2: NUM_BIRDS_DOA_PLANT -> "50", but stored as text (matching the rest of
# the sheet, which keeps numeric-looking values as text) and without adding
# any new cell style. A plain Value assignment would auto-coerce "50" into a
# real number, so instead stage the text in a scratch cell (quote-prefixed
# so it's unambiguously text), copy it, and paste-special *values only* into
# L2 so the text type carries over but the scratch cell's formatting does not.
$scratch = $ws.Range("E2")
$scratch.Value = "'50"
$scratch.Copy()
$ws.Range("L2").PasteSpecial(-4163)
$scratch.Clear()

# AM2: BIRD_SIZE -> "Pullet" (plain text, no numeric ambiguity)
$ws.Range("AM2").Value = "Pullet"
